$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.654.91"
$ws.Range("E2").Value = "  +2.75%  "
$ws.Range("D3").Value = "3.601.63"
$ws.Range("E3").Value = "  +4.95%  "
$ws.Range("E4").Value = "  -0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "238.73"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.44%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "658.87"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +6.34%  "
$ws.Range("E7").Value = "  +7.76%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.407"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +3.70%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  +4.74%  "
$ws.Range("D11").Value = "3.600.13"
$ws.Range("E11").Value = "  +4.89%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "43.12"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("E13").Value = "  +0.95%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.32"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("D15").Value = "4.289.00"
$ws.Range("E15").Value = "  +5.39%  "
$ws.Range("D16").Value = "95.468.56"
$ws.Range("E16").Value = "  +2.68%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.0000255"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +3.90%  "
$ws.Range("D18").Value = "3.603.67"
$ws.Range("E18").Value = "  +5.14%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.80"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -4.12%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.64"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +8.80%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "18.04"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.22%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "3.62"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +8.34%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.496"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +11.79%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "510.47"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.23%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.0000196"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +6.23%  "
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("E27").Value = "  +6.82%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "12.72"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +6.71%  "
$ws.Range("D29").Value = "3.778.36"
$ws.Range("E29").Value = "  +4.49%  "
$ws.Range("E30").Value = "  +15.54%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "11.33"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.68%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.00%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.140"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.46%  "
$ws.Range("E34").Value = "  -0.37%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.177"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.60%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "31.86"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +5.45%  "
$ws.Range("E37").Value = "  +2.74%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "576.39"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.87%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "8.17"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +9.99%  "
$ws.Range("E40").Value = "  +6.65%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.151"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.01%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.924"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("E44").Value = "  +4.44%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.73"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.31%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "23.78"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.46%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "33.80"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +30.44%  "
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("E49").Value = "  +6.52%  "
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("E51").Value = "  -6.09%  "
